$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Building data")

# Update the header label in A1 to reflect the renamed field used for
# mapping/merging (jurisdiction_taxlot_identifier -> jurisdiction_tax_lot_id)
$ws.Range("A1").Value = "jurisdiction_tax_lot_id"
